# Leave Card update (1/5/2024 4:46 PM)
# - Adds a "SL(3-0-0)" entry for the 12/31/2023 pay period (row 39)
# - Starts a new "2024" year section (row 40)
# - Extends the PERIOD date column through the existing blank rows (41-86)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 40: new "2024" year-section header (matches the look of the ---
# --- existing "2022" / "2023" headers in A10 / A23)                  ---
$ws.Range("A23").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A40").Value = "'2024"

# --- Row 39: new Sick Leave entry for the 12/31/2023 period ---
$ws.Range("B39").Value = "SL(3-0-0)"
$ws.Range("C39").Value = 1.25
$ws.Range("K39").Value = "12/28, 29,31/2023"

# --- Rows 41-86: fill in the (previously blank) PERIOD end-of-month dates ---
$periodDates = @{
    41 = 45322; 42 = 45351; 43 = 45382; 44 = 45412; 45 = 45443; 46 = 45473;
    47 = 45504; 48 = 45535; 49 = 45565; 50 = 45596; 51 = 45626; 52 = 45657;
    53 = 45688; 54 = 45716; 55 = 45747; 56 = 45777; 57 = 45808; 58 = 45838;
    59 = 45869; 60 = 45900; 61 = 45930; 62 = 45961; 63 = 45991; 64 = 46022;
    65 = 46053; 66 = 46081; 67 = 46112; 68 = 46142; 69 = 46173; 70 = 46203;
    71 = 46234; 72 = 46265; 73 = 46295; 74 = 46326; 75 = 46356; 76 = 46387;
    77 = 46418; 78 = 46446; 79 = 46477; 80 = 46507; 81 = 46538; 82 = 46568;
    83 = 46599; 84 = 46630; 85 = 46660; 86 = 46691
}

foreach ($row in $periodDates.Keys) {
    $ws.Range("A$row").Value = $periodDates[$row]
}

$excel.CalculateFull()
